$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0004425917103382563
$ws.Range("C2").Value = 0.0004624487404803749
$ws.Range("D2").Value = 0.0004922456954451093
$ws.Range("E2").Value = 0.0004920847268673355
$ws.Range("F2").Value = 0.0004723011363636364
$ws.Range("G2").Value = 0.0004087928464977645
$ws.Range("H2").Value = 0.0003410049627791564

$ws.Range("B3").Value = 0.007224868985231061
$ws.Range("C3").Value = 0.007457527826596369
$ws.Range("D3").Value = 0.007833679325925021
$ws.Range("E3").Value = 0.007785953177257524
$ws.Range("F3").Value = 0.007361505681818182
$ws.Range("G3").Value = 0.006696472925981123
$ws.Range("H3").Value = 0.005490694789081887

$ws.Range("B4").Value = 0.1143401619818961
$ws.Range("C4").Value = 0.120796719390744
$ws.Range("D4").Value = 0.1267554035901819
$ws.Range("E4").Value = 0.1266220735785953
$ws.Range("F4").Value = 0.1198863636363637
$ws.Range("G4").Value = 0.1061102831594635
$ws.Range("H4").Value = 0.08740694789081886

$ws.Range("B5").Value = 0.5724154359218675
$ws.Range("C5").Value = 0.6115992970123023
$ws.Range("D5").Value = 0.6335327878861887
$ws.Range("E5").Value = 0.6494983277591972
$ws.Range("F5").Value = 0.6100852272727273
$ws.Range("G5").Value = 0.5375062096373572
$ws.Range("H5").Value = 0.4362282878411911

$ws.Range("B7").Value = 1.829204383039542
$ws.Range("C7").Value = 1.932630345635618
$ws.Range("D7").Value = 2.008792282329955
$ws.Range("E7").Value = 2.018729096989966
$ws.Range("F7").Value = 1.905894886363636
$ws.Range("G7").Value = 1.692498758072528
$ws.Range("H7").Value = 1.385235732009926

$ws.Range("B8").Value = 4.454502143878037
$ws.Range("C8").Value = 4.699472759226714
$ws.Range("D8").Value = 4.952985712541214
$ws.Range("E8").Value = 4.958751393534002
$ws.Range("F8").Value = 4.655539772727272
$ws.Range("G8").Value = 4.18628912071535
$ws.Range("H8").Value = 3.406327543424318

$ws.Range("B9").Value = 9.163887565507384
$ws.Range("C9").Value = 9.648506151142357
$ws.Range("D9").Value = 10.28697032604714
$ws.Range("E9").Value = 10.33444816053512
$ws.Range("F9").Value = 9.7265625
$ws.Range("G9").Value = 8.628912071535023
$ws.Range("H9").Value = 7.096774193548388

$ws.Range("B10").Value = 17.09623630300143
$ws.Range("C10").Value = 17.69771528998243
$ws.Range("D10").Value = 18.91561851263891
$ws.Range("E10").Value = 18.93199554069119
$ws.Range("F10").Value = 17.90838068181818
$ws.Range("G10").Value = 15.73770491803279
$ws.Range("H10").Value = 13.13275434243176
